# Update cryptocurrency price (D) and 1h-volume-change (E) columns
# per the scraper refresh commit "Updated cryptos list on Mon Jun 24 23:42:57 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.297.74"
$ws.Range("E2").Value = "  -4.77%  "
$ws.Range("D3").Value = "3.347.83"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'566.53"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "'132.28"
$ws.Range("E6").Value = "  +3.04%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.342.95"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "'7.46"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  -2.54%  "
$ws.Range("D12").Value = "'0.378"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "3.914.16"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "3.348.29"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "'24.85"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "60.347.83"
$ws.Range("E18").Value = "  -4.85%  "
$ws.Range("D19").Value = "'13.61"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "'5.72"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "'9.27"
$ws.Range("E21").Value = "  -5.38%  "
$ws.Range("D22").Value = "'363.44"
$ws.Range("E22").Value = "  -5.03%  "
$ws.Range("D23").Value = "'0.561"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "3.474.30"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D26").Value = "'69.70"
$ws.Range("E26").Value = "  -5.28%  "
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("E28").Value = "  +18.17%  "
$ws.Range("D29").Value = "'7.54"
$ws.Range("E29").Value = "  +7.87%  "
$ws.Range("D30").Value = "'0.987"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").Value = "'8.03"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").Value = "'2.14"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "3.374.09"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").Value = "'22.98"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("D38").Value = "'6.93"
$ws.Range("E38").Value = "  +3.18%  "
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "'158.88"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("D41").Value = "'0.0776"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "'4.40"
$ws.Range("E43").Value = "  +2.57%  "
$ws.Range("E44").Value = "  +10.23%  "
$ws.Range("D45").Value = "'40.95"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "'0.751"
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("D47").Value = "'23.62"
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "'6.83"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").Value = "'22.65"
$ws.Range("E50").Value = "  +12.35%  "
$ws.Range("D51").Value = "'0.898"
$ws.Range("E51").Value = "  +2.16%  "
